$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.628.64"
$ws.Range("E2").Value = "  +1.41%  "

$ws.Range("D3").Value = "'1.801.12"
$ws.Range("E3").Value = "  +1.34%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").Value = "'227.53"
$ws.Range("E5").Value = "  +0.80%  "

$ws.Range("E6").Value = "  +2.01%  "

$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("E8").Value = "  +4.12%  "

$ws.Range("E9").Value = "  +1.89%  "

$ws.Range("E10").Value = "  +0.92%  "

$ws.Range("E11").Value = "  +0.45%  "

$ws.Range("D12").Value = "'2.063.26"
$ws.Range("E12").Value = "  +1.42%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.816.21"
$ws.Range("E13").Value = "  +2.47%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'11.18"
$ws.Range("E14").Value = "  +2.77%  "

$ws.Range("E15").Value = "  +3.24%  "

$ws.Range("D16").Value = "'34.601.40"
$ws.Range("E16").Value = "  +1.43%  "

$ws.Range("E17").Value = "  +3.81%  "

$ws.Range("D18").Value = "'69.01"
$ws.Range("E18").Value = "  +1.79%  "

$ws.Range("E19").Value = "  +0.76%  "

$ws.Range("D20").Value = "'247.42"
$ws.Range("E20").Value = "  +0.91%  "

$ws.Range("D21").Value = "'11.37"
$ws.Range("E21").Value = "  +3.83%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.32%  "

$ws.Range("E23").Value = "  +2.46%  "

$ws.Range("D24").Value = "'171.58"
$ws.Range("E24").Value = "  +5.64%  "

$ws.Range("E25").Value = "  +1.59%  "

$ws.Range("E26").Value = "  +2.14%  "

$ws.Range("D27").Value = "'16.64"
$ws.Range("E27").Value = "  +2.28%  "

$ws.Range("E28").Value = "  +2.20%  "

$ws.Range("E29").Value = "  -0.22%  "

$ws.Range("D30").Value = "'4.08"
$ws.Range("E30").Value = "  +10.51%  "

$ws.Range("E31").Value = "  +1.66%  "

$ws.Range("E32").Value = "  +0.70%  "

$ws.Range("E33").Value = "  +2.48%  "

$ws.Range("E34").Value = "  +2.76%  "

$ws.Range("D35").Value = "'1.434.81"
$ws.Range("E35").Value = "  -0.14%  "

$ws.Range("D36").Value = "'2.57"
$ws.Range("E36").Value = "  +7.34%  "

$ws.Range("D37").Value = "'0.676"
$ws.Range("E37").Value = "  +2.62%  "

$ws.Range("E38").Value = "  +2.61%  "

$ws.Range("E39").Value = "  +0.75%  "

$ws.Range("D40").Value = "'84.95"
$ws.Range("E40").Value = "  +6.12%  "

$ws.Range("D41").Value = "'0.948"
$ws.Range("E41").Value = "  +3.23%  "

$ws.Range("E42").Value = "  +1.46%  "

$ws.Range("D43").Value = "'2.75"
$ws.Range("E43").Value = "  +3.25%  "

$ws.Range("D44").Value = "'13.79"
$ws.Range("E44").Value = "  +2.48%  "

$ws.Range("E45").Value = "  +3.10%  "

$ws.Range("E46").Value = "  +0.77%  "

$ws.Range("E47").Value = "  +0.57%  "

$ws.Range("D48").Value = "'1.963.95"
$ws.Range("E48").Value = "  +1.38%  "

$ws.Range("D49").Value = "'105.42"
$ws.Range("E49").Value = "  +1.25%  "

$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  -0.24%  "

$ws.Range("E51").Value = "  -5.17%  "
